# Update "Estado de Cuenta" workbook:
#  - Remove the 3 other workers' data (Yesenia, Cristel, Sandy)
#  - Keep only Keisy Sierra Rincon, now covering periods 1702..2012 (47 periods)
#  - Update summary header cells (Valor Mora total, worker count, period count)
#  - Remove now-unused trailing rows, shifting the footer notes rows up

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Propagate the special "last data row" border styling (currently on
#        row 71) onto row 62, which will become the new last data row once
#        the extra rows are deleted below. ---
$ws.Range("B71:J71").Copy() | Out-Null
$ws.Range("B62:J62").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- 2. Delete the now-superfluous rows 63-71 (old filler + old last row).
#        This shifts the trailing "firma" notes rows (76,77) up to (67,68). ---
$ws.Range("A63:A71").EntireRow.Delete() | Out-Null

# --- 3. Header / summary cells ---
$ws.Range("E11").Value = 1722576   # VALOR MORA total
$ws.Range("C13").Value = 1         # Cant. Trabajadores
$ws.Range("F13").Value = 47        # Cant. Periodos

# --- 4. Rewrite the data table (rows 16-62) for Keisy Sierra Rincon across
#        periods 1702 .. 2012, ascending. ---
$periods = @("1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712", `
             "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812", `
             "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912", `
             "2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012")

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = $startRow + $i
    $period = $periods[$i]

    $ws.Cells.Item($r, 2).Value = "CC"                     # B: Tipo Doc Trabajador
    $ws.Cells.Item($r, 3).Value = "1085226925"             # C: N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = "KEISY SIERRA RINCON"    # D: Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $period                  # E: Periodo Mora

    if ($r -eq ($startRow + $periods.Length - 1)) {
        $ws.Cells.Item($r, 6).Value = 25820                # F: Valor Mora (last period)
    } else {
        $ws.Cells.Item($r, 6).Value = 36886                # F: Valor Mora
    }
    $ws.Cells.Item($r, 7).Value = 922133                   # G: Salario Basico
}
